$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values are written with a leading apostrophe to
# force Excel to store them as literal text (prices/percentages would
# otherwise be auto-converted to numbers), then the style is reset to
# 'Normal' so the quote-prefix flag / number format isn't left behind.
$updates = @(
    @{cell='D2'; value='28.093.13'},
    @{cell='E2'; value='  -1.52%  '},
    @{cell='D3'; value='1.894.77'},
    @{cell='E3'; value='  -0.82%  '},
    @{cell='E4'; value='  -0.03%  '},
    @{cell='D5'; value='314.74'},
    @{cell='E5'; value='  +0.13%  '},
    @{cell='D6'; value='1.001'},
    @{cell='E6'; value='  +0.00%  '},
    @{cell='D7'; value='0.5017'},
    @{cell='E7'; value='  -0.69%  '},
    @{cell='D8'; value='0.3900'},
    @{cell='E8'; value='  -1.33%  '},
    @{cell='D9'; value='0.09213'},
    @{cell='E9'; value='  -5.62%  '},
    @{cell='D10'; value='1.128'},
    @{cell='E10'; value='  -2.71%  '},
    @{cell='D11'; value='41.89'},
    @{cell='E11'; value='  +0.19%  '},
    @{cell='D12'; value='6.379'},
    @{cell='E12'; value='  -2.59%  '},
    @{cell='D13'; value='20.79'},
    @{cell='E13'; value='  -1.77%  '},
    @{cell='D14'; value='1.902.81'},
    @{cell='E14'; value='  -0.28%  '},
    @{cell='D15'; value='7.278'},
    @{cell='E16'; value='  -0.02%  '},
    @{cell='B17'; value='Litecoin'},
    @{cell='C17'; value='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{cell='D17'; value='92.45'},
    @{cell='E17'; value='  -1.50%  '},
    @{cell='B18'; value='ShibaInu'},
    @{cell='C18'; value='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'},
    @{cell='D18'; value='0.00001107'},
    @{cell='E18'; value='  -2.77%  '},
    @{cell='D19'; value='0.06651'},
    @{cell='E19'; value='  -0.07%  '},
    @{cell='D20'; value='17.83'},
    @{cell='E20'; value='  -1.59%  '},
    @{cell='D21'; value='1.001'},
    @{cell='E21'; value='  +0.00%  '},
    @{cell='D22'; value='6.204'},
    @{cell='E22'; value='  -1.53%  '},
    @{cell='D23'; value='28.155.39'},
    @{cell='E23'; value='  -1.54%  '},
    @{cell='D24'; value='11.44'},
    @{cell='E24'; value='  -0.14%  '},
    @{cell='D25'; value='2.319'},
    @{cell='E25'; value='  +1.74%  '},
    @{cell='B26'; value='WrappedliquidstakedEther2.0'},
    @{cell='C26'; value='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'},
    @{cell='D26'; value='2.125.85'},
    @{cell='E26'; value='  -0.10%  '},
    @{cell='B27'; value='LidoDAOToken'},
    @{cell='C27'; value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'},
    @{cell='D27'; value='2.551'},
    @{cell='E27'; value='  -7.54%  '},
    @{cell='B28'; value='Monero'},
    @{cell='C28'; value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{cell='D28'; value='158.30'},
    @{cell='E28'; value='  -0.55%  '},
    @{cell='B29'; value='EthereumClassic'},
    @{cell='C29'; value='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'},
    @{cell='D29'; value='20.83'},
    @{cell='E29'; value='  -2.07%  '},
    @{cell='B30'; value='BitcoinCash'},
    @{cell='C30'; value='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'},
    @{cell='D30'; value='126.98'},
    @{cell='E30'; value='  -1.41%  '},
    @{cell='B31'; value='ImmutableX'},
    @{cell='C31'; value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{cell='D31'; value='1.078'},
    @{cell='E31'; value='  -2.19%  '},
    @{cell='B32'; value='Stellar'},
    @{cell='C32'; value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{cell='D32'; value='0.1056'},
    @{cell='E32'; value='  -1.40%  '},
    @{cell='B33'; value='Filecoin'},
    @{cell='C33'; value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{cell='D33'; value='5.605'},
    @{cell='E33'; value='  -2.34%  '},
    @{cell='B34'; value='HuobiToken'},
    @{cell='C34'; value='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'},
    @{cell='D34'; value='3.615'},
    @{cell='E34'; value='  -0.71%  '},
    @{cell='B35'; value='FraxShare'},
    @{cell='C35'; value='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'},
    @{cell='D35'; value='9.540'},
    @{cell='E35'; value='  -3.12%  '},
    @{cell='B36'; value='TrustWalletToken'},
    @{cell='C36'; value='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'},
    @{cell='D36'; value='1.347'},
    @{cell='E36'; value='  +13.23%  '},
    @{cell='B37'; value='Hedera'},
    @{cell='C37'; value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{cell='D37'; value='0.06594'},
    @{cell='E37'; value='  -2.70%  '},
    @{cell='B38'; value='VeChain'},
    @{cell='C38'; value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{cell='D38'; value='0.02400'},
    @{cell='E38'; value='  -1.76%  '},
    @{cell='B39'; value='Algorand'},
    @{cell='C39'; value='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'},
    @{cell='D39'; value='0.2203'},
    @{cell='E39'; value='  -1.20%  '},
    @{cell='B40'; value='ARBITRUM'},
    @{cell='C40'; value='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'},
    @{cell='D40'; value='1.222'},
    @{cell='E40'; value='  -3.76%  '},
    @{cell='B41'; value='TheSandbox'},
    @{cell='C41'; value='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{cell='D41'; value='0.6458'},
    @{cell='E41'; value='  +0.38%  '},
    @{cell='B42'; value='Aptos'},
    @{cell='C42'; value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{cell='D42'; value='11.42'},
    @{cell='E42'; value='  -2.56%  '},
    @{cell='B43'; value='InternetComputer(DFINITY)'},
    @{cell='C43'; value='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'},
    @{cell='D43'; value='4.970'},
    @{cell='E43'; value='  -2.40%  '},
    @{cell='B44'; value='Frax'},
    @{cell='C44'; value='https://coinranking.com/coin/KfWtaeV1W+frax-frax'},
    @{cell='D44'; value='1.001'},
    @{cell='E44'; value='  -0.02%  '},
    @{cell='B45'; value='Decentraland'},
    @{cell='C45'; value='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'},
    @{cell='D45'; value='0.6097'},
    @{cell='E45'; value='  +0.19%  '},
    @{cell='B46'; value='EnergySwap'},
    @{cell='C46'; value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{cell='D46'; value='13.29'},
    @{cell='E46'; value='  -3.27%  '},
    @{cell='B47'; value='WEMIXTOKEN'},
    @{cell='C47'; value='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'},
    @{cell='D47'; value='1.307'},
    @{cell='E47'; value='  +2.25%  '},
    @{cell='B48'; value='PancakeSwap'},
    @{cell='C48'; value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
    @{cell='D48'; value='3.693'},
    @{cell='E48'; value='  +0.60%  '},
    @{cell='B49'; value='NEARProtocol'},
    @{cell='C49'; value='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'},
    @{cell='D49'; value='1.999'},
    @{cell='E49'; value='  -2.12%  '},
    @{cell='B50'; value='Quant'},
    @{cell='C50'; value='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'},
    @{cell='D50'; value='122.22'},
    @{cell='E50'; value='  -2.01%  '},
    @{cell='B51'; value='EOS'},
    @{cell='C51'; value='https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'},
    @{cell='D51'; value='1.202'},
    @{cell='E51'; value='  -0.91%  '}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.cell)
    $c.Value = "'" + $u.value
    $c.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells"
